# Corecciones lógicas y gráficas
# - Reorder / correct the "Materias primas" recipe strings (Harina now listed
#   first, with the quantities corrected; Torta's Leche goes from 5.0 to 3.0).
# - Add a new "Disponible" column (E) filled with 1 for every product row.
# - Select cell D7 to mirror the saved selection in the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("E1").Value = "Disponible"

# Corrected "Materias primas" text per row (column C) -------------------------
$ws.Range("C2").Value = "2.0-Harina  (kg),1.0-Vainilla (ml),3.0-Leche (litros),1.0-Huevos (unidad),"
$ws.Range("C3").Value = "2.0-Harina  (kg),5.0-Huevos (unidad),"
$ws.Range("C4").Value = "5.0-Harina  (kg),1.0-Vainilla (ml),2.0-Huevos (unidad),"
$ws.Range("C5").Value = "5.0-Harina  (kg),5.0-Huevos (unidad),"
$ws.Range("C6").Value = "4.0-Harina  (kg),2.0-Limon (unidad),5.0-Crema (litros),5.0-Huevos (unidad),"
$ws.Range("C7").Value = "0.3-Harina  (kg),0.2-Leche (litros),0.1-Vainilla (ml),2.0-Huevos (unidad),"

# New "Disponible" column (E) -------------------------------------------------
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 1

# Match saved selection in the source workbook
$ws.Range("D7").Select()
